# Generate Report for Handoff
# Marks the 8c254ff2 file row as "Ready for handoff" across Overview, zh-cn, de-de sheets,
# records the new handoff datetime, and records the stale-handback error detail.

$wb = $excel.ActiveWorkbook

# --- Overview sheet (row 3 = 8c254ff2-a83e-465d-b8e4-238f5c507734.md) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-20 00:52:42"

# --- zh-cn sheet (row 3 = 8c254ff2-a83e-465d-b8e4-238f5c507734.md) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-20 00:52:37"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/594a203fa0a07d476ed2c499042acea1f5d7c28d/e2e/8c254ff2-a83e-465d-b8e4-238f5c507734.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d44ea3820eb08f406431c0e26044d0d276a969a1/e2e/8c254ff2-a83e-465d-b8e4-238f5c507734.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666667

# --- de-de sheet (row 3 = 8c254ff2-a83e-465d-b8e4-238f5c507734.md) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-20 00:52:42"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/594a203fa0a07d476ed2c499042acea1f5d7c28d/e2e/8c254ff2-a83e-465d-b8e4-238f5c507734.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d44ea3820eb08f406431c0e26044d0d276a969a1/e2e/8c254ff2-a83e-465d-b8e4-238f5c507734.md."
$dede.Columns.Item(16).ColumnWidth = 39.16666666666667
